# Replace the arithmetic problems in the practice-sheet table with the new
# set of problems, preserving cell formatting (font/size) by only touching
# the run text via Find/Replace within each cell's range.

$d = $word.ActiveDocument

$newValues = @(
    "114÷9=",
    "965÷7=",
    "135÷5=",
    "231÷7=",
    "325÷3=",
    "855÷5=",
    "470÷3=",
    "397÷2=",
    "661÷8=",
    "102÷6=",
    "106÷9=",
    "988÷4=",
    "681÷7=",
    "196÷5=",
    "478÷9=",
    "177÷5=",
    "686÷9=",
    "592÷2=",
    "626÷2=",
    "288÷3=",
    "373÷5=",
    "899÷8=",
    "271÷3=",
    "955÷4=",
    "987÷8="
)

$table = $d.Tables.Item(1)
$idx = 0

foreach ($row in $table.Rows) {
    foreach ($cell in $row.Cells) {
        $text = $cell.Range.Text
        # Cell text ends with the cell-mark/bell control chars; trim those
        # to detect whether the cell actually contains a problem string.
        $trimmed = $text.TrimEnd([char]7, [char]13)
        if ($trimmed -ne "") {
            $newText = $newValues[$idx]
            $cellRange = $cell.Range
            $cellRange.MoveEnd(1, -1) | Out-Null
            $cellRange.Text = $newText
            $idx = $idx + 1
        }
    }
}
